# "backing up all files and data" -- rename the qtrs/region tabs to their
# fuller backup-friendly names, move the active selection on each sheet,
# and refresh the bestFit-ish width on the month/quarter-number columns.

$wb = $excel.ActiveWorkbook

# --- rename worksheets -------------------------------------------------
$wsTraveled = $wb.Worksheets.Item("total_miles_traveled")
$wsQtrs     = $wb.Worksheets.Item("qtrs")
$wsRegion   = $wb.Worksheets.Item("region")

$wsQtrs.Name   = "total_miles_quarters"
$wsRegion.Name = "total_miles_region"

# --- column width touch-up ---------------------------------------------
# total_miles_traveled: give the Month column (B) a snug, best-fit-like width
$wsTraveled.Columns.Item(2).ColumnWidth = 5.92

# total_miles_quarters: widen qtr_num column (C) now that it shows its header
$wsQtrs.Columns.Item(3).ColumnWidth = 7.59

# --- selections / active tab --------------------------------------------
# Move the cursor on the two non-active sheets first ...
$wsQtrs.Activate()
$wsQtrs.Range("F8").Select()

$wsRegion.Activate()
$wsRegion.Range("G33").Select()

# ... then land on total_miles_traveled, which becomes the active/selected tab
$wsTraveled.Activate()
$wsTraveled.Range("H12").Select()
